$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (BI_SYMBOL / PD_SYMBOL / UM_SYMBOL / MP_SYMBOL) - left to right
$ws.Range("AC2").Value = "BI001"
$ws.Range("AD2").Value = "PD001"
$ws.Range("AE2").Value = "UM001"
$ws.Range("AF2").Value = "MP001"

# Row 3 - left to right
$ws.Range("AC3").Value = "BI002"
$ws.Range("AD3").Value = "PD002"
$ws.Range("AE3").Value = "UM002"
$ws.Range("AF3").Value = "MP002"

# Rows 4 & 5, column by column
$ws.Range("AC4").Value = "BI003"
$ws.Range("AC5").Value = "BI004"

$ws.Range("AD4").Value = "PD003"
$ws.Range("AD5").Value = "PD004"

$ws.Range("AE4").Value = "UM003"
$ws.Range("AE5").Value = "UM004"

$ws.Range("AF4").Value = "MP003"
$ws.Range("AF5").Value = "MP004"

# VERSION column (B2:B5) updated last
$ws.Range("B2:B5").Value = "SYMBOL_2017"

# Update the selection to match the new view state
$ws.Range("E9").Select()
